$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item("Patient_day")

# ----------------------------------------------------------------------
# New "report parameter" header block (F1:G4) on the Patient_day sheet.
# Formats are cloned from an existing bordered/cambria cell (B11) so the
# shared font/border definitions in styles.xml are reused rather than
# duplicated, then tweaked (fill / number format) to match the new look.
# ----------------------------------------------------------------------

# Yellow header cells F1 (report) / G1 (parameter)
$ws2.Range("B11").Copy()
$ws2.Range("F1:G1").PasteSpecial(-4122)
$ws2.Range("F1:G1").Interior.Color = 65535

# Plain bordered label/value cells
$ws2.Range("B11").Copy()
$ws2.Range("F2").PasteSpecial(-4122)
$ws2.Range("F3").PasteSpecial(-4122)
$ws2.Range("F4").PasteSpecial(-4122)
$ws2.Range("G4").PasteSpecial(-4122)

# Date-formatted value cells G2 (start date) / G3 (end date)
$ws2.Range("B11").Copy()
$ws2.Range("G2:G3").PasteSpecial(-4122)
$ws2.Range("G2:G3").NumberFormat = "dd\-mm\-yy"

$excel.CutCopyMode = $false

# Cell contents -- written in this particular order so the shared-string
# table picks up new unique strings in the same sequence as the original
# authoring session (Start date, End date, parameter, Hospital, report,
# Siem Reap ... Hospital).
$ws2.Range("F2").Value = "Start date"
$ws2.Range("G2").Value = [DateTime]"2019-01-01"
$ws2.Range("F3").Value = "End date"
$ws2.Range("G3").Value = [DateTime]"2019-12-31"
$ws2.Range("G1").Value = "parameter"
$ws2.Range("F4").Value = "Hospital"
$ws2.Range("F1").Value = "report"
$ws2.Range("G4").Value = "Siem Reap Provincial Referral  Hospital"

# Column widths for the new columns
$ws2.Columns.Item(6).ColumnWidth = 13.166666666666666
$ws2.Columns.Item(7).ColumnWidth = 35

# Input-message-only data validations on G2 / G3 (no restriction type,
# just a prompt + "information" style error/input message), matching the
# order they were authored in (G3 first, then G2).
$dv3 = $ws2.Range("G3").Validation
$dv3.Delete()
$dv3.Add(0, 3, 1, [Type]::Missing, [Type]::Missing)
$dv3.IgnoreBlank = $true
$dv3.InputTitle = "information"
$dv3.InputMessage = "dd-mm-yy"
$dv3.ErrorMessage = "dd-mm-yy"
$dv3.ShowInput = $true
$dv3.ShowError = $true

$dv2 = $ws2.Range("G2").Validation
$dv2.Delete()
$dv2.Add(0, 3, 1, [Type]::Missing, [Type]::Missing)
$dv2.IgnoreBlank = $true
$dv2.InputTitle = "information"
$dv2.InputMessage = "dd-mm-yy"
$dv2.ShowInput = $true
$dv2.ShowError = $true

# Make Patient_day the active/selected sheet & selected cell (G6), and
# clear the tab-selected flag that used to live on the Dict sheet.
$ws2.Activate()
$ws2.Range("G6").Select()
